$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing parameter values in row 2
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 70
$ws.Range("F2").Value = 0.0005
$ws.Range("H2").Value = 10

# Add new transformer parameter headers
$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"

# Add new transformer parameter values
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 1

# Update selection to match target view state
$ws.Range("J4").Select()
